$d = $word.ActiveDocument

# 1. Remove the trailing standalone space run after the "Tip: ..." sentence.
$d.Content.Find.Execute("fødselsår. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "fødselsår.", 2)

# 2. Remove the word "din" from ", så din andre kan se på din bruger, hvor du kommer fra."
$d.Content.Find.Execute("så din andre", $true, $false, $false, $false, $false,
                         $true, 1, $false, "så andre", 2)
